$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.650.56"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.444.62"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'570.19"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'144.93"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").Value = "2.440.23"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "'5.24"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").Value = "'27.12"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  -4.34%  "
$ws.Range("D16").Value = "2.885.14"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "62.558.95"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "2.431.13"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").Value = "'11.24"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "'327.28"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  +11.63%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'65.14"
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").Value = "'629.30"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "'9.11"
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("D28").Value = "0.0₃0999"
$ws.Range("E28").Value = "  -4.83%  "
$ws.Range("D29").Value = "2.558.65"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = "  -4.54%  "
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "'5.13"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").Value = "'1.51"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").Value = "'146.58"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("D43").Value = "'2.59"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D45").Value = "'146.65"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").Value = "'3.76"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'20.75"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").Value = "'0.597"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "'0.0922"
$ws.Range("E51").Value = "  -0.65%  "
